# "Generate Report for Handoff" — stamp the latest handoff-xliff generation
# timestamp for the file that was just (re)handed off
# (6bd30ad5-963c-4f40-a8b9-900428c2e088.md) across the Overview sheet and
# each per-locale status sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 6bd30ad5... row (row 6).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-08-19 12:41:35"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the 6bd30ad5... row
# (row 6).
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-08-19 12:41:31"

# de-de sheet: "Latest Handoff Datetime" column (H) for the 6bd30ad5... row
# (row 6).
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-08-19 12:41:35"
